# Generate Report for Handoff
# Replace the GUID-named handoff artifact (old GUID -> new GUID) and refresh
# the associated handoff/handback timestamps across the Overview, zh-cn and
# de-de sheets, including the "display" text of the hyperlinks that point at
# the markdown file on GitHub.

$wb = $excel.ActiveWorkbook

$oldGuidFile   = "1d7502aa-c67a-4212-88de-37eb5e6c16c2"
$newGuidFile   = "f1f4970f-6358-4806-9c46-bf5b6792414c"

$oldZhXlf = "$oldGuidFile.13253dd8b9cb193cab726e1ce0eff5fbe5a3bc01.zh-cn.xlf"
$newZhXlf = "$newGuidFile.36ee70e9debff942dff1769005c135e5ef2db414.zh-cn.xlf"

$oldDeXlf = "$oldGuidFile.13253dd8b9cb193cab726e1ce0eff5fbe5a3bc01.de-de.xlf"
$newDeXlf = "$newGuidFile.36ee70e9debff942dff1769005c135e5ef2db414.de-de.xlf"

$githubBase = "https://github.com/OpenLocalizationTestOrg/ol-test1/blob/24bb51156d274aca6e720ac3e125d3c6fb41d088/e2e/"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("A2").Value = "$newGuidFile.md"
$overview.Range("B2").Value = "e2e\$newGuidFile.md"
$overview.Range("G2").Value = "2017-01-03 06:21:33"

$overview.Range("B2").Hyperlinks.Delete()
$overview.Hyperlinks.Add($overview.Range("B2"), ($githubBase + "$newGuidFile.md"), "", "", "e2e\$newGuidFile.md")

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("A2").Value = "$newGuidFile.md"
$zhcn.Range("G2").Value = $newZhXlf
$zhcn.Range("H2").Value = "2017-01-03 06:21:23"

$zhcn.Range("A2").Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), ($githubBase + "$newGuidFile.md"), "", "", "$newGuidFile.md")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("A2").Value = "$newGuidFile.md"
$dede.Range("G2").Value = $newDeXlf
$dede.Range("H2").Value = "2017-01-03 06:21:33"

$dede.Range("A2").Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), ($githubBase + "$newGuidFile.md"), "", "", "$newGuidFile.md")

Write-Output "Handoff report regenerated: $oldGuidFile -> $newGuidFile"
